$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 11, column B ("Rule" label for the R40 row) is updated from "R40" to "1".
# Use a leading apostrophe so Excel stores this numeric-looking entry as TEXT
# (shared string "1") instead of converting it to a Number value.
$cell = $ws.Range("B11")
$cell.Value = "'1"
